# Word COM-interop script implementing the diff.
$d = $word.ActiveDocument

# 1) Merge the three runs around "homo sapiens" into a single run whose
#    text is the quoted phrase with no intervening run breaks.
$d.Content.Find.Execute("“homo sapiens“", $true, $false, $false, $false, $false,
                         $true, 1, $false, "“homo sapiens“", 2) | Out-Null

# 2) "October 2014)" -> "January 2015)" (first occurrence, narrative text).
$d.Content.Find.Execute("October 2014)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "January 2015)", 2) | Out-Null

# 3) "october-2014.fasta" -> "january-2015.fasta" (first filename mention).
$d.Content.Find.Execute("october-2014.fasta", $true, $false, $false, $false, $false,
                         $true, 1, $false, "january-2015.fasta", 2) | Out-Null

# 4) "uniprot-human-reviewed-september-2014.fasta" -> "uniprot-human-reviewed-january-2014.fasta"
$d.Content.Find.Execute("uniprot-human-reviewed-september-2014.fasta", $true, $false, $false, $false, $false,
                         $true, 1, $false, "uniprot-human-reviewed-january-2014.fasta", 2) | Out-Null

# 5) "october-2014.fasta" -> "january-2015.fasta" (second filename mention, trypsin file).
$d.Content.Find.Execute("trypsin-october-2014.fasta", $true, $false, $false, $false, $false,
                         $true, 1, $false, "trypsin-january-2015.fasta", 2) | Out-Null
